$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# --- Data corrections, rows 204-208 (column P only) ---
$ws.Range("P204").Value = 330
$ws.Range("P205").Value = 339
$ws.Range("P206").Value = 379
$ws.Range("P207").Value = 428
$ws.Range("P208").Value = 438

# --- row 209 ---
$ws.Range("N209").Value = 114
$ws.Range("P209").Value = 460

# --- row 210 ---
$ws.Range("N210").Value = 119
$ws.Range("P210").Value = 447

# --- row 211 ---
$ws.Range("N211").Value = 117
$ws.Range("O211").Value = 404
$ws.Range("P211").Value = 438

# --- row 212 ---
$ws.Range("N212").Value = 110
$ws.Range("O212").Value = 427
$ws.Range("P212").Value = 474

# --- row 213 ---
$ws.Range("N213").Value = 108
$ws.Range("O213").Value = 445
$ws.Range("P213").Value = 479

# --- row 214 ---
$ws.Range("N214").Value = 96
$ws.Range("O214").Value = 390
$ws.Range("P214").Value = 515

# --- row 215 ---
$ws.Range("N215").Value = 95
$ws.Range("O215").Value = 346
$ws.Range("P215").Value = 531

# --- row 216 ---
$ws.Range("C216").Value = 20
$ws.Range("F216").Value = 1
$ws.Range("G216").Value = 7
$ws.Range("O216").Value = 307
$ws.Range("P216").Value = 548

# --- row 217 ---
$ws.Range("C217").Value = 7
$ws.Range("F217").Value = 1
$ws.Range("G217").Value = 7
$ws.Range("I217").Value = 1
$ws.Range("N217").Value = 88
$ws.Range("O217").Value = 301
$ws.Range("P217").Value = 491

# --- row 218: previously blank (only the date + "" formulas), now filled in ---
$ws.Range("C218").Value = 0
$ws.Range("D218").Value = 0
$ws.Range("E218").Value = 1
$ws.Range("F218").Value = 1
$ws.Range("G218").Value = 7
$ws.Range("I218").Value = 0
$ws.Range("L218").Value = "0"
$ws.Range("M218").Value = "0"
$ws.Range("N218").Value = 76
$ws.Range("O218").Value = 221
$ws.Range("P218").Value = 412

# --- sheet view: move the active selection of the frozen (bottom-right) pane ---
$ws.Range("R164").Select()
